$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Insert()

$links = @(
    "https://www.conservationlaos.com/",
    "https://saeda.net/"
)
$ws.Range("A1").Hyperlinks.Delete()
for ($i = 0; $i -lt $links.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $ws.Hyperlinks.Add($cell, $links[$i])
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.Underline = $true
    $cell.Font.Color = 16711680
}
Write-Host $ws.Cells.Item(2,1).Font.Name
Write-Host $ws.Cells.Item(2,1).Font.Underline
Write-Host $ws.Cells.Item(2,1).Font.Color
